$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns retain text formatting so values
# such as "25.717.60" or "14.83" are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.717.60"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").Value = "1.745.94"
$ws.Range("E3").Value = "  -5.22%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "238.68"
$ws.Range("E5").Value = "  -8.34%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4938"
$ws.Range("E7").Value = "  -6.35%  "
$ws.Range("D8").Value = "41.52"
$ws.Range("E8").Value = "  -7.60%  "
$ws.Range("D9").Value = "0.2472"
$ws.Range("E9").Value = "  -21.89%  "
$ws.Range("D10").Value = "0.05972"
$ws.Range("E10").Value = "  -12.16%  "
$ws.Range("D11").Value = "1.743.83"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D12").Value = "0.06792"
$ws.Range("E12").Value = "  -12.91%  "
$ws.Range("D13").Value = "14.83"
$ws.Range("D14").Value = "4.469"
$ws.Range("E14").Value = "  -10.86%  "
$ws.Range("D15").Value = "77.21"
$ws.Range("E15").Value = "  -12.56%  "
$ws.Range("D16").Value = "0.5830"
$ws.Range("E16").Value = "  -25.60%  "
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "25.758.35"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -16.94%  "
$ws.Range("D21").Value = "0.000006508"
$ws.Range("E21").Value = "  -17.96%  "
$ws.Range("D22").Value = "1.967.78"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").Value = "3.976"
$ws.Range("E23").Value = "  -13.80%  "
$ws.Range("D24").Value = "7.908"
$ws.Range("E24").Value = "  -15.28%  "
$ws.Range("D25").Value = "5.026"
$ws.Range("E25").Value = "  -16.18%  "
$ws.Range("D26").Value = "136.17"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").Value = "1.491"
$ws.Range("E27").Value = "  -12.84%  "
$ws.Range("D28").Value = "1.835"
$ws.Range("E28").Value = "  -17.35%  "
$ws.Range("D29").Value = "14.56"
$ws.Range("E29").Value = "  -14.44%  "
$ws.Range("D30").Value = "100.96"
$ws.Range("E30").Value = "  -8.90%  "
$ws.Range("D31").Value = "3.804"
$ws.Range("E31").Value = "  -9.74%  "
$ws.Range("D32").Value = "0.08103"
$ws.Range("E32").Value = "  -6.85%  "
$ws.Range("D33").Value = "3.352"
$ws.Range("E33").Value = "  -17.85%  "
$ws.Range("D34").Value = "0.04415"
$ws.Range("E34").Value = "  -9.19%  "
$ws.Range("D35").Value = "0.9986"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "2.647"
$ws.Range("E36").Value = "  -7.65%  "
$ws.Range("D37").Value = "1.018"
$ws.Range("E37").Value = "  -10.55%  "
$ws.Range("D38").Value = "0.6062"
$ws.Range("E38").Value = "  -17.19%  "
$ws.Range("D39").Value = "2.701"
$ws.Range("D40").Value = "2.058"
$ws.Range("E40").Value = "  -12.21%  "
$ws.Range("D42").Value = "103.55"
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("D43").Value = "0.01499"
$ws.Range("E43").Value = "  -13.45%  "
$ws.Range("D44").Value = "0.7789"
$ws.Range("E44").Value = "  -13.87%  "
$ws.Range("D45").Value = "5.187"
$ws.Range("E45").Value = "  -12.27%  "
$ws.Range("D46").Value = "0.3768"
$ws.Range("E46").Value = "  -21.90%  "
$ws.Range("D47").Value = "0.05125"
$ws.Range("E47").Value = "  -11.98%  "
$ws.Range("D48").Value = "0.1081"
$ws.Range("E48").Value = "  -13.20%  "
$ws.Range("D49").Value = "5.960"
$ws.Range("E49").Value = "  -22.67%  "
$ws.Range("D50").Value = "30.31"
$ws.Range("E50").Value = "  -13.17%  "
$ws.Range("E51").Value = "  -12.50%  "
